# Applies the fohbot/leng.xlsx "test" sheet edits:
#  - A10 text changes ("Здесь вы можете настроить свои каналы и создать
#    публикации" -> "Всего каналов: ")
#  - A12 text changes (old "Отправьте инфо для поста..." short prompt ->
#    new, much longer post-content instructions) and its row grows from
#    60pt to 405pt tall to fit the new text
#  - the sheet's active selection moves from C33 to C12
#
# (NB: the shared-string table reflow seen in the target OOXML diff -
# strings 9/11 dropped, the two new strings appended at the end, and every
# other row's <v> shifting to match - is a side effect Excel performs
# automatically when strings are replaced; we only need to set the two
# cell values below and it falls out on save.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Всего каналов: "

$postText = @'
Отправьте содержимое поста (Текст, фото, видео, аудио, документ, gif, хеш-тег)
• Используйте * (звездочку)  в начале и в конце текста чтобы сделать <b>жирный шрифт</b>.
• Используйте _ (нижнее подчеркивание)  в начале и в конце текста чтобы сделать <i>курсивный шрифт</i>.
• Используйте ` (апостроф) в начале и в конце текста чтобы сделать <code>выделенный текст</code>.
• Для того чтобы использовать Ссылку в тексте напишите в любом месте поста в скобках [Наш канал](https://t.me/Fohbot_News).<a href="https://t.me/Fohbot_News">Наш канал</a>
'@
$ws.Range("A12").Value = $postText
$ws.Rows.Item(12).RowHeight = 405

# Move the selection the way the author's session ended up (C33 -> C12).
[void]$ws.Range("C12").Select()
